$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pod")
$ws.Range("D13").Value = 999
Write-Host "done"
